# Update sigma_010 (sheet2) and sigma_025 (sheet3) with refined values, and
# add a new sigma_050 sheet with its own noisy/denoised PSNR data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# sigma_010 : refine columns B (Noisy) and C (NLM-LBP) for rows 2-12
# ---------------------------------------------------------------------------
$ws010 = $wb.Worksheets.Item("sigma_010")

$sigma010 = @{
    2  = @(27.64330152387554, 30.18409227405678)
    3  = @(27.6388194136901,  30.14679154490338)
    4  = @(27.65746032602065, 30.17788033786534)
    5  = @(27.61541982907814, 30.19340504797643)
    6  = @(27.6133364115809,  30.17548416185557)
    7  = @(27.6111291900843,  30.18267121769808)
    8  = @(27.63400777090692, 30.16462014897921)
    9  = @(27.65835848296324, 30.17191001769045)
    10 = @(27.61881785038445, 30.16487577435396)
    11 = @(27.63769825191305, 30.17422750189417)
    12 = @(27.63283490504973, 30.17359580272734)
}

foreach ($row in $sigma010.Keys) {
    $vals = $sigma010[$row]
    $ws010.Cells.Item($row, 2).Value2 = $vals[0]
    $ws010.Cells.Item($row, 3).Value2 = $vals[1]
}

# ---------------------------------------------------------------------------
# sigma_025 : refine columns B (Noisy) and C (NLM-LBP) for rows 2-12
# ---------------------------------------------------------------------------
$ws025 = $wb.Worksheets.Item("sigma_025")

$sigma025 = @{
    2  = @(19.71954964720373, 26.74317763887622)
    3  = @(19.72047349950375, 26.71779814047028)
    4  = @(19.72700060860124, 26.71905242006973)
    5  = @(19.74216498042597, 26.72443608683247)
    6  = @(19.74778316055646, 26.76175223935482)
    7  = @(19.74467392601867, 26.72791028781838)
    8  = @(19.71769075644981, 26.7101891381126)
    9  = @(19.7475168699232,  26.72877213239066)
    10 = @(19.73329653627579, 26.77725931346652)
    11 = @(19.73767604810264, 26.74962482809111)
    12 = @(19.73378260330612, 26.73599722254828)
}

foreach ($row in $sigma025.Keys) {
    $vals = $sigma025[$row]
    $ws025.Cells.Item($row, 2).Value2 = $vals[0]
    $ws025.Cells.Item($row, 3).Value2 = $vals[1]
}

# ---------------------------------------------------------------------------
# sigma_050 : new sheet appended at the end, same layout as the others
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws050 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws050.Name = "sigma_050"

$ws050.Cells.Item(1, 1).Value = "Rows"
$ws050.Cells.Item(1, 2).Value = "Noisy"
$ws050.Cells.Item(1, 3).Value = "NLM-LBP"

$sigma050 = @{
    2  = @(0,       14.84612988285806, 21.41504394480652)
    3  = @(1,       14.81461750231181, 21.46399440157916)
    4  = @(2,       14.8344562335508,  21.4486173816412)
    5  = @(3,       14.8358983617873,  21.48508168049317)
    6  = @(4,       14.8227928762877,  21.43749951796696)
    7  = @(5,       14.82265100569516, 21.43158318782671)
    8  = @(6,       14.82030694092514, 21.45498149451899)
    9  = @(7,       14.83351771380503, 21.48882317938185)
    10 = @(8,       14.8295499738282,  21.48983256651341)
    11 = @(9,       14.83102723393659, 21.48930892665396)
}

foreach ($row in $sigma050.Keys) {
    $vals = $sigma050[$row]
    $ws050.Cells.Item($row, 1).Value2 = $vals[0]
    $ws050.Cells.Item($row, 2).Value2 = $vals[1]
    $ws050.Cells.Item($row, 3).Value2 = $vals[2]
}

$ws050.Cells.Item(12, 1).Value = "Média"
$ws050.Cells.Item(12, 2).Value2 = 14.82909477249858
$ws050.Cells.Item(12, 3).Value2 = 21.46047662813819

$ws010.Select()
Write-Output "done"
